$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1031.72
$ws.Range("I28").Value = 985.86664
$ws.Range("K28").Value = 985.86664
$ws.Range("M28").Value = -500.86664

$ws.Range("H62").Value = 2366.6667
$ws.Range("I62").Value = 2612.5
$ws.Range("J62").Value = 1875
$ws.Range("K62").Value = 2612.5
$ws.Range("L62").Value = 1875
$ws.Range("M62").Value = -1988.5
$ws.Range("N62").Value = -3123

$ws.Range("H65").Value = 2366.6667
$ws.Range("I65").Value = 2612.5
$ws.Range("J65").Value = 1875
$ws.Range("K65").Value = 13062.5
$ws.Range("L65").Value = 9375
$ws.Range("M65").Value = -9942.5
$ws.Range("N65").Value = -15615

$ws.Range("H138").Value = 2637.8333
$ws.Range("I138").Value = 3030.2
$ws.Range("J138").Value = 2398.5854
$ws.Range("K138").Value = 9090.599999999999
$ws.Range("L138").Value = 7195.7562
$ws.Range("M138").Value = -3950.599999999999
$ws.Range("N138").Value = -17475.7562

$ws.Range("H140").Value = 35000
$ws.Range("I140").Value = 35000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 35000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -29820
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 2801.1428
$ws.Range("I141").Value = 1127
$ws.Range("J141").Value = 5033.3335
$ws.Range("K141").Value = 3381
$ws.Range("L141").Value = 15100.0005
$ws.Range("M141").Value = 1799
$ws.Range("N141").Value = -25460.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3135.6
$ws.Range("I32").Value = 3016.9387
$ws.Range("J32").Value = 8950
$ws.Range("K32").Value = 3016.9387
$ws.Range("L32").Value = 8950
$ws.Range("M32").Value = -2729.9387
$ws.Range("N32").Value = -9524

$ws.Range("H44").Value = 30250
$ws.Range("J44").Value = 30250
$ws.Range("L44").Value = 30250
$ws.Range("N44").Value = -31226

$ws.Range("H110").Value = 1620.0714
$ws.Range("I110").Value = 1610.2
$ws.Range("J110").Value = 1644.75
$ws.Range("K110").Value = 1610.2
$ws.Range("L110").Value = 1644.75
$ws.Range("M110").Value = 434.8
$ws.Range("N110").Value = -5734.75

$ws.Range("H122").Value = 1594
$ws.Range("I122").Value = 1502.8334
$ws.Range("J122").Value = 1703.4
$ws.Range("K122").Value = 4508.5002
$ws.Range("L122").Value = 5110.200000000001
$ws.Range("M122").Value = -2058.5002
$ws.Range("N122").Value = -10010.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 500
$ws.Range("I94").Value = 500
$ws.Range("K94").Value = 500
$ws.Range("M94").Value = -49

$ws.Range("H105").Value = 1775.7142
$ws.Range("I105").Value = 1664.4445
$ws.Range("J105").Value = 1976
$ws.Range("K105").Value = 1664.4445
$ws.Range("L105").Value = 1976
$ws.Range("M105").Value = 82.55549999999994
$ws.Range("N105").Value = -5470

$ws.Range("H107").Value = 4471.6924
$ws.Range("I107").Value = 3567.3333
$ws.Range("J107").Value = 6506.5
$ws.Range("K107").Value = 3567.3333
$ws.Range("L107").Value = 6506.5
$ws.Range("M107").Value = -1647.3333
$ws.Range("N107").Value = -10346.5

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws.Range("H141").Value = 57911.125
$ws.Range("J141").Value = 51430
$ws.Range("L141").Value = 51430
$ws.Range("N141").Value = -61790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1443.375
$ws.Range("I16").Value = 3100
$ws.Range("J16").Value = 891.1667
$ws.Range("K16").Value = 3100
$ws.Range("L16").Value = 891.1667
$ws.Range("M16").Value = -2813
$ws.Range("N16").Value = -1465.1667

$ws.Range("H31").Value = 2282.5217
$ws.Range("I31").Value = 1338.9333
$ws.Range("K31").Value = 1338.9333
$ws.Range("M31").Value = -1043.9333

$ws.Range("H34").Value = 2282.5217
$ws.Range("I34").Value = 1338.9333
$ws.Range("K34").Value = 1338.9333
$ws.Range("M34").Value = -1136.9333

$ws.Range("H51").Value = 9010.875
$ws.Range("J51").Value = 10142.429
$ws.Range("L51").Value = 10142.429
$ws.Range("N51").Value = -11614.429

$ws.Range("H58").Value = 2427.84
$ws.Range("I58").Value = 1209.8667
$ws.Range("J58").Value = 4254.8
$ws.Range("K58").Value = 1209.8667
$ws.Range("L58").Value = 4254.8
$ws.Range("M58").Value = -1006.8667
$ws.Range("N58").Value = -4660.8

$ws.Range("H61").Value = 9010.875
$ws.Range("J61").Value = 10142.429
$ws.Range("L61").Value = 10142.429
$ws.Range("N61").Value = -10838.429

$ws.Range("H74").Value = 15288.875
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 16901.572
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 16901.572
$ws.Range("M74").Value = -3126
$ws.Range("N74").Value = -18649.572

$ws.Range("H77").Value = 15288.875
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 16901.572
$ws.Range("K77").Value = 12000
$ws.Range("L77").Value = 50704.716
$ws.Range("M77").Value = -7632
$ws.Range("N77").Value = -59440.716

$ws.Range("H99").Value = 73074.71000000001
$ws.Range("I99").Value = 37353
$ws.Range("J99").Value = 144518.14
$ws.Range("K99").Value = 37353
$ws.Range("L99").Value = 144518.14
$ws.Range("M99").Value = -35855
$ws.Range("N99").Value = -147514.14

$ws.Range("H105").Value = 797.9
$ws.Range("I105").Value = 797.9
$ws.Range("K105").Value = 797.9
$ws.Range("M105").Value = 949.1

$ws.Range("H113").Value = 1443.375
$ws.Range("I113").Value = 3100
$ws.Range("J113").Value = 891.1667
$ws.Range("K113").Value = 3100
$ws.Range("L113").Value = 891.1667
$ws.Range("M113").Value = -930
$ws.Range("N113").Value = -5231.1667

$ws.Range("H122").Value = 4812866
$ws.Range("I122").Value = 7360057.5
$ws.Range("J122").Value = 1504.4445
$ws.Range("K122").Value = 22080172.5
$ws.Range("L122").Value = 4513.333500000001
$ws.Range("M122").Value = -22077722.5
$ws.Range("N122").Value = -9413.333500000001

$ws.Range("H126").Value = 73074.71000000001
$ws.Range("I126").Value = 37353
$ws.Range("J126").Value = 144518.14
$ws.Range("K126").Value = 112059
$ws.Range("L126").Value = 433554.42
$ws.Range("M126").Value = -109589
$ws.Range("N126").Value = -438494.42

$ws.Range("H136").Value = 2427.84
$ws.Range("I136").Value = 1209.8667
$ws.Range("J136").Value = 4254.8
$ws.Range("K136").Value = 3629.6001
$ws.Range("L136").Value = 12764.4
$ws.Range("M136").Value = -1079.6001
$ws.Range("N136").Value = -17864.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 835.01514
$ws.Range("J5").Value = 982.90247
$ws.Range("L5").Value = 2948.70741
$ws.Range("N5").Value = -3172.70741

$ws.Range("H103").Value = 3778260.2
$ws.Range("I103").Value = 4250505
$ws.Range("J103").Value = 300
$ws.Range("K103").Value = 12751515
$ws.Range("L103").Value = 900
$ws.Range("M103").Value = -12750636
$ws.Range("N103").Value = -2658

$ws.Range("H107").Value = 117963.12
$ws.Range("I107").Value = 91154.45
$ws.Range("J107").Value = 167112.33
$ws.Range("K107").Value = 273463.35
$ws.Range("L107").Value = 501336.99
$ws.Range("M107").Value = -271543.35
$ws.Range("N107").Value = -505176.99

$ws.Range("H131").Value = 2926.6345
$ws.Range("J131").Value = 1843.1
$ws.Range("L131").Value = 5529.299999999999
$ws.Range("N131").Value = -15609.3

$ws.Range("H135").Value = 835.01514
$ws.Range("J135").Value = 982.90247
$ws.Range("L135").Value = 8846.122230000001
$ws.Range("N135").Value = -13916.12223

$ws.Range("H140").Value = 1728.55
$ws.Range("I140").Value = 1182.3846
$ws.Range("J140").Value = 2742.8572
$ws.Range("K140").Value = 3547.1538
$ws.Range("L140").Value = 8228.571599999999
$ws.Range("M140").Value = 1632.8462
$ws.Range("N140").Value = -18588.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1175.5454
$ws.Range("I113").Value = 1103.4445
$ws.Range("K113").Value = 1103.4445
$ws.Range("M113").Value = 1066.5555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 409.69232
$ws.Range("I22").Value = 507.5
$ws.Range("K22").Value = 507.5
$ws.Range("M22").Value = -212.5

$ws.Range("H27").Value = 409.69232
$ws.Range("I27").Value = 507.5
$ws.Range("K27").Value = 507.5
$ws.Range("M27").Value = -400.5

$ws.Range("H40").Value = 2361.889
$ws.Range("I40").Value = 2361.889
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2361.889
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2225.889
$ws.Range("N40").ClearContents()

$ws.Range("H45").Value = 12987.4
$ws.Range("I45").Value = 10970.5
$ws.Range("J45").Value = 14332
$ws.Range("K45").Value = 10970.5
$ws.Range("L45").Value = 14332
$ws.Range("M45").Value = -10563.5
$ws.Range("N45").Value = -15146

$ws.Range("H48").Value = 12499
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 10000
$ws.Range("M48").Value = -9339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 14998
$ws.Range("J50").Value = 14998
$ws.Range("L50").Value = 14998
$ws.Range("N50").Value = -16260

$ws.Range("H54").Value = 9755.111000000001
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 16759.2
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 16759.2
$ws.Range("M54").Value = -480
$ws.Range("N54").Value = -17799.2

$ws.Range("H132").Value = 4365.48
$ws.Range("I132").Value = 5419.457
$ws.Range("K132").Value = 16258.371
$ws.Range("M132").Value = -13728.371
